# Update the "Förändrad" (Changed) date column (C) for every data row
# from 2026-02-07 (serial 46060) to 2026-02-08 (serial 46061).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = [DateTime]::FromOADate(46061)

for ($row = 2; $row -le 509; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46060) {
        $cell.Value = $newDate
    }
}
